# Fixed a bug in linesTrigger.
#
# The rows in the data sheet were paired with the wrong symbol/line
# ids (e.g. row 2 held symbol 1203's reel data instead of symbol
# 401's). Re-point every data row (2-21) at its correct symbol id and
# reel1..reel5 counts so each row's data matches its symbol again.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# rowNumber -> [symbol, reel1, reel2, reel3, reel4, reel5]
$rows = @(2, 3, 4, 5, 6, 7, 8, 9, 10, 11, 12, 13, 14, 15, 16, 17, 18, 19, 20, 21)
$data = @(
    @(401, 9, 48, 67, 75, 45),
    @(601, 9, 60, 67, 60, 42),
    @(201, 9, 30, 15, 45, 30),
    @(1202, 2, 10, 10, 10, 10),
    @(101, 9, 30, 15, 60, 15),
    @(501, 9, 52, 30, 75, 45),
    @(1201, 2, 10, 10, 10, 10),
    @(301, 6, 45, 30, 60, 45),
    @(701, 3, 90, 45, 97, 15),
    @(1001, 18, 30, 75, 60, 72),
    @(902, 1, 0, 0, 0, 0),
    @(801, 3, 67, 65, 52, 45),
    @(1203, 3, 15, 15, 15, 15),
    @(901, 16, 15, 45, 60, 60),
    @(1, 0, 2, 2, 2, 2),
    @(1101, 0, 15, 30, 30, 0),
    @(2, 0, 2, 2, 2, 2),
    @(502, 0, 4, 0, 0, 0),
    @(802, 0, 4, 5, 4, 0),
    @(3, 0, 3, 3, 3, 3)
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $rowNum = $rows[$i]
    $vals = $data[$i]
    for ($j = 0; $j -lt $vals.Length; $j++) {
        $col = [char](65 + $j)
        $ws.Range("$col$rowNum").Value = $vals[$j]
    }
}
